$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 80.5
$ws.Range("I11").Value = 80.5
$ws.Range("K11").Value = 80.5
$ws.Range("M11").Value = 59.5

$ws.Range("H15").Value = 1885.5714
$ws.Range("I15").Value = 1885.5714
$ws.Range("K15").Value = 5656.7142
$ws.Range("M15").Value = -5487.7142

$ws.Range("H28").Value = 375
$ws.Range("I28").Value = 375
$ws.Range("K28").Value = 375
$ws.Range("M28").Value = 110

$ws.Range("H88").Value = 1749.8572
$ws.Range("I88").Value = 2197.5
$ws.Range("J88").Value = 1570.8
$ws.Range("K88").Value = 2197.5
$ws.Range("L88").Value = 1570.8
$ws.Range("M88").Value = -1791.5
$ws.Range("N88").Value = -2382.8

$ws.Range("H91").Value = 1749.8572
$ws.Range("I91").Value = 2197.5
$ws.Range("J91").Value = 1570.8
$ws.Range("K91").Value = 2197.5
$ws.Range("L91").Value = 1570.8
$ws.Range("M91").Value = -793.5
$ws.Range("N91").Value = -4378.8

$ws.Range("H103").Value = 2095.4
$ws.Range("I103").Value = 1879
$ws.Range("J103").Value = 2149.5
$ws.Range("K103").Value = 5637
$ws.Range("L103").Value = 6448.5
$ws.Range("M103").Value = -5051
$ws.Range("N103").Value = -7620.5

$ws.Range("H106").Value = 2371.6667
$ws.Range("I106").Value = 2371.6667
$ws.Range("K106").Value = 2371.6667
$ws.Range("M106").Value = -1740.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9513.643
$ws.Range("I32").Value = 7501.825
$ws.Range("K32").Value = 7501.825
$ws.Range("M32").Value = -7214.825

$ws.Range("H45").Value = 2946
$ws.Range("I45").Value = 2946
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 2946
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -2569
$ws.Range("N45").ClearContents()

$ws.Range("H63").Value = 4954.3335
$ws.Range("I63").Value = 705.6667
$ws.Range("K63").Value = 705.6667
$ws.Range("M63").Value = -19.66669999999999

$ws.Range("H66").Value = 4954.3335
$ws.Range("I66").Value = 705.6667
$ws.Range("K66").Value = 3528.3335
$ws.Range("M66").Value = -96.33349999999973

$ws.Range("H74").Value = 14391.444
$ws.Range("I74").Value = 13400.6
$ws.Range("K74").Value = 13400.6
$ws.Range("M74").Value = -12526.6

$ws.Range("H77").Value = 14391.444
$ws.Range("I77").Value = 13400.6
$ws.Range("K77").Value = 67003
$ws.Range("M77").Value = -62635

$ws.Range("H88").Value = 1976.1111
$ws.Range("I88").Value = 896.25
$ws.Range("J88").Value = 2840
$ws.Range("K88").Value = 896.25
$ws.Range("L88").Value = 2840
$ws.Range("M88").Value = -490.25
$ws.Range("N88").Value = -3652

$ws.Range("H91").Value = 1976.1111
$ws.Range("I91").Value = 896.25
$ws.Range("J91").Value = 2840
$ws.Range("K91").Value = 896.25
$ws.Range("L91").Value = 2840
$ws.Range("M91").Value = 507.75
$ws.Range("N91").Value = -5648

$ws.Range("H110").Value = 3296.0667
$ws.Range("I110").Value = 1402.6364
$ws.Range("K110").Value = 1402.6364
$ws.Range("M110").Value = 642.3635999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 35265.332
$ws.Range("J82").Value = 75283
$ws.Range("L82").Value = 75283
$ws.Range("N82").Value = -76049

$ws.Range("H85").Value = 35265.332
$ws.Range("J85").Value = 75283
$ws.Range("L85").Value = 75283
$ws.Range("N85").Value = -77935

$ws.Range("H94").Value = 3630.0557
$ws.Range("I94").Value = 3488.6428
$ws.Range("K94").Value = 3488.6428
$ws.Range("M94").Value = -3037.6428

$ws.Range("H108").Value = 59666.668
$ws.Range("J108").Value = 59666.668
$ws.Range("L108").Value = 59666.668
$ws.Range("N108").Value = -67346.66800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2176.5
$ws.Range("I31").Value = 1412.625
$ws.Range("K31").Value = 1412.625
$ws.Range("M31").Value = -1117.625

$ws.Range("H34").Value = 2176.5
$ws.Range("I34").Value = 1412.625
$ws.Range("K34").Value = 1412.625
$ws.Range("M34").Value = -1210.625

$ws.Range("H58").Value = 4134.3
$ws.Range("I58").Value = 2852.1667
$ws.Range("K58").Value = 2852.1667
$ws.Range("M58").Value = -2649.1667

$ws.Range("H99").Value = 4217.857
$ws.Range("I99").Value = 4269.6665
$ws.Range("K99").Value = 4269.6665
$ws.Range("M99").Value = -2771.6665

$ws.Range("H126").Value = 4217.857
$ws.Range("I126").Value = 4269.6665
$ws.Range("K126").Value = 12808.9995
$ws.Range("M126").Value = -10338.9995

$ws.Range("H132").Value = 3592.125
$ws.Range("I132").Value = 3010.75
$ws.Range("K132").Value = 9032.25
$ws.Range("M132").Value = -6502.25

$ws.Range("H136").Value = 4134.3
$ws.Range("I136").Value = 2852.1667
$ws.Range("K136").Value = 8556.500100000001
$ws.Range("M136").Value = -6006.500100000001

$ws.Range("H141").Value = 618791.8
$ws.Range("I141").Value = 300000
$ws.Range("J141").Value = 698489.75
$ws.Range("K141").Value = 300000
$ws.Range("L141").Value = 698489.75
$ws.Range("M141").Value = -294820
$ws.Range("N141").Value = -708849.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 122.833336
$ws.Range("I38").Value = 129.4
$ws.Range("J38").Value = 90
$ws.Range("K38").Value = 388.2
$ws.Range("L38").Value = 270
$ws.Range("M38").Value = -41.20000000000005
$ws.Range("N38").Value = -964

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

$ws.Range("H70").Value = 3999.6667
$ws.Range("I70").Value = 999.5
$ws.Range("K70").Value = 999.5
$ws.Range("M70").Value = -729.5

$ws.Range("H73").Value = 3999.6667
$ws.Range("I73").Value = 999.5
$ws.Range("K73").Value = 999.5
$ws.Range("M73").Value = -63.5

$ws.Range("H132").Value = 3998.8333
$ws.Range("I132").Value = 2999
$ws.Range("J132").Value = 4998.6665
$ws.Range("K132").Value = 8997
$ws.Range("L132").Value = 14995.9995
$ws.Range("M132").Value = -6467
$ws.Range("N132").Value = -20055.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2222
$ws.Range("J46").Value = 2222
$ws.Range("L46").Value = 2222
$ws.Range("N46").Value = -2598

$ws.Range("H132").Value = 3669.5
$ws.Range("I132").Value = 2893
$ws.Range("K132").Value = 8679
$ws.Range("M132").Value = -6149

$ws.Range("H136").Value = 3955.3076
$ws.Range("I136").Value = 1440.5555
$ws.Range("K136").Value = 4321.666499999999
$ws.Range("M136").Value = -1771.666499999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 852.7778
$ws.Range("I136").Value = 848.375
$ws.Range("K136").Value = 2545.125
$ws.Range("M136").Value = 4.875
